# Auto-generated edit script
# 1) Refresh market-price derived columns (H..N) per row across all 8 profession sheets
# 2) Strip the bold/bordered/centered header style from row 1 on every sheet (revert to default style)

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1431.6666
$ws.Range("I18").Value = 1431.6666
$ws.Range("K18").Value = 1431.6666
$ws.Range("M18").Value = -1147.6666
$ws.Range("H61").Value = 3146.8235
$ws.Range("I61").Value = 1750
$ws.Range("J61").Value = 3333.0667
$ws.Range("K61").Value = 5250
$ws.Range("L61").Value = 9999.2001
$ws.Range("M61").Value = -5078
$ws.Range("N61").Value = -10343.2001
$ws.Range("H70").Value = 5490.4375
$ws.Range("I70").Value = 5549.5
$ws.Range("K70").Value = 16648.5
$ws.Range("M70").Value = -16378.5
$ws.Range("H73").Value = 5490.4375
$ws.Range("I73").Value = 5549.5
$ws.Range("K73").Value = 16648.5
$ws.Range("M73").Value = -15712.5
$ws.Range("H82").Value = 1306.3334
$ws.Range("I82").Value = 1012.2632
$ws.Range("J82").Value = 4100
$ws.Range("K82").Value = 3036.7896
$ws.Range("L82").Value = 12300
$ws.Range("M82").Value = -2630.7896
$ws.Range("N82").Value = -13112
$ws.Range("H85").Value = 1306.3334
$ws.Range("I85").Value = 1012.2632
$ws.Range("J85").Value = 4100
$ws.Range("K85").Value = 3036.7896
$ws.Range("L85").Value = 12300
$ws.Range("M85").Value = -1632.7896
$ws.Range("N85").Value = -15108
$ws.Range("H86").Value = 4291.1763
$ws.Range("J86").Value = 4812.1816
$ws.Range("L86").Value = 4812.1816
$ws.Range("N86").Value = -7058.1816
$ws.Range("H89").Value = 4291.1763
$ws.Range("J89").Value = 4812.1816
$ws.Range("L89").Value = 24060.908
$ws.Range("N89").Value = -35292.908
$ws.Range("H98").Value = 3293.7856
$ws.Range("I98").Value = 3915.7
$ws.Range("J98").Value = 1739
$ws.Range("K98").Value = 3915.7
$ws.Range("L98").Value = 1739
$ws.Range("M98").Value = -2417.7
$ws.Range("N98").Value = -4735
$ws.Range("H99").Value = 436
$ws.Range("I99").Value = 375.33334
$ws.Range("K99").Value = 1126.00002
$ws.Range("M99").Value = 371.9999800000001
$ws.Range("H100").Value = 7599.5264
$ws.Range("I100").Value = 5878.9
$ws.Range("J100").Value = 9511.333
$ws.Range("K100").Value = 5878.9
$ws.Range("L100").Value = 9511.333
$ws.Range("M100").Value = -5337.9
$ws.Range("N100").Value = -10593.333
$ws.Range("H101").Value = 932.6667
$ws.Range("I101").Value = 882.3333
$ws.Range("J101").Value = 983
$ws.Range("K101").Value = 2646.9999
$ws.Range("L101").Value = 2949
$ws.Range("M101").Value = -1024.9999
$ws.Range("N101").Value = -6193
$ws.Range("H113").Value = 9880
$ws.Range("J113").Value = 12919
$ws.Range("L113").Value = 12919
$ws.Range("N113").Value = -19427
$ws.Range("H122").Value = 3293.7856
$ws.Range("I122").Value = 3915.7
$ws.Range("J122").Value = 1739
$ws.Range("K122").Value = 11747.1
$ws.Range("L122").Value = 5217
$ws.Range("M122").Value = -9297.099999999999
$ws.Range("N122").Value = -10117
$ws.Range("H137").Value = 11125923
$ws.Range("I137").Value = 20025018
$ws.Range("K137").Value = 60075054
$ws.Range("M137").Value = -60072504
$ws.Range("H138").Value = 4629.8535
$ws.Range("J138").Value = 4887.971
$ws.Range("L138").Value = 14663.913
$ws.Range("N138").Value = -24943.913

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15389165
$ws.Range("I32").Value = 16671478
$ws.Range("J32").Value = 1408.6
$ws.Range("K32").Value = 16671478
$ws.Range("L32").Value = 1408.6
$ws.Range("M32").Value = -16671191
$ws.Range("N32").Value = -1982.6
$ws.Range("H45").Value = 4502.857
$ws.Range("I45").Value = 3694.1052
$ws.Range("J45").Value = 6210.222
$ws.Range("K45").Value = 3694.1052
$ws.Range("L45").Value = 6210.222
$ws.Range("M45").Value = -3317.1052
$ws.Range("N45").Value = -6964.222
$ws.Range("H61").Value = 3426.25
$ws.Range("I61").Value = 3426.25
$ws.Range("K61").Value = 3426.25
$ws.Range("M61").Value = -3214.25
$ws.Range("H122").Value = 4883.909
$ws.Range("I122").Value = 4604.375
$ws.Range("J122").Value = 5147
$ws.Range("K122").Value = 13813.125
$ws.Range("L122").Value = 15441
$ws.Range("M122").Value = -11363.125
$ws.Range("N122").Value = -20341
$ws.Range("H132").Value = 2498
$ws.Range("I132").Value = 2004.5
$ws.Range("K132").Value = 6013.5
$ws.Range("M132").Value = -3483.5
$ws.Range("H136").Value = 3426.25
$ws.Range("I136").Value = 3426.25
$ws.Range("K136").Value = 10278.75
$ws.Range("M136").Value = -7728.75

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 51809.8
$ws.Range("I20").Value = 1783.4166
$ws.Range("K20").Value = 1783.4166
$ws.Range("M20").Value = -1536.4166
$ws.Range("H94").Value = 3699.6924
$ws.Range("I94").Value = 3631.889
$ws.Range("J94").Value = 3852.25
$ws.Range("K94").Value = 3631.889
$ws.Range("L94").Value = 3852.25
$ws.Range("M94").Value = -3180.889
$ws.Range("N94").Value = -4754.25
$ws.Range("H105").Value = 8248.25
$ws.Range("I105").Value = 7996.4
$ws.Range("J105").Value = 8668
$ws.Range("K105").Value = 7996.4
$ws.Range("L105").Value = 8668
$ws.Range("M105").Value = -6249.4
$ws.Range("N105").Value = -12162
$ws.Range("H107").Value = 7676.375
$ws.Range("J107").Value = 7000
$ws.Range("L107").Value = 7000
$ws.Range("N107").Value = -10840
$ws.Range("H117").Value = 90000
$ws.Range("J117").Value = 90000
$ws.Range("L117").Value = 90000
$ws.Range("N117").Value = -99178

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 1312.25
$ws.Range("I5").Value = 1273.9
$ws.Range("K5").Value = 1273.9
$ws.Range("M5").Value = -1161.9
$ws.Range("H7").Value = 152.12
$ws.Range("I7").Value = 94.23529
$ws.Range("K7").Value = 94.23529
$ws.Range("M7").Value = 18.76470999999999
$ws.Range("H16").Value = 11799
$ws.Range("I16").Value = 13998.75
$ws.Range("K16").Value = 13998.75
$ws.Range("M16").Value = -13711.75
$ws.Range("H22").Value = 347
$ws.Range("I22").Value = 347
$ws.Range("K22").Value = 347
$ws.Range("M22").Value = 3
$ws.Range("H31").Value = 6333.32
$ws.Range("I31").Value = 11385.8
$ws.Range("K31").Value = 11385.8
$ws.Range("M31").Value = -11090.8
$ws.Range("H34").Value = 6333.32
$ws.Range("I34").Value = 11385.8
$ws.Range("K34").Value = 11385.8
$ws.Range("M34").Value = -11183.8
$ws.Range("H58").Value = 2747.5625
$ws.Range("I58").Value = 2423
$ws.Range("K58").Value = 2423
$ws.Range("M58").Value = -2220
$ws.Range("H113").Value = 11799
$ws.Range("I113").Value = 13998.75
$ws.Range("K113").Value = 13998.75
$ws.Range("M113").Value = -11828.75
$ws.Range("H136").Value = 2747.5625
$ws.Range("I136").Value = 2423
$ws.Range("K136").Value = 7269
$ws.Range("M136").Value = -4719

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 26359296
$ws.Range("I4").Value = 30235506
$ws.Range("J4").Value = 1075.8
$ws.Range("K4").Value = 90706518
$ws.Range("L4").Value = 3227.4
$ws.Range("M4").Value = -90706406
$ws.Range("N4").Value = -3451.4
$ws.Range("H36").Value = 718
$ws.Range("I36").Value = 218.33333
$ws.Range("J36").Value = 1217.6666
$ws.Range("K36").Value = 654.99999
$ws.Range("L36").Value = 3652.9998
$ws.Range("M36").Value = -485.99999
$ws.Range("N36").Value = -3990.9998
$ws.Range("H68").Value = 1266.8975
$ws.Range("I68").Value = 1079.6666
$ws.Range("J68").Value = 1282.5
$ws.Range("K68").Value = 3238.9998
$ws.Range("L68").Value = 3847.5
$ws.Range("M68").Value = -2427.9998
$ws.Range("N68").Value = -5469.5
$ws.Range("H71").Value = 1266.8975
$ws.Range("I71").Value = 1079.6666
$ws.Range("J71").Value = 1282.5
$ws.Range("K71").Value = 9716.9994
$ws.Range("L71").Value = 11542.5
$ws.Range("M71").Value = -5660.999400000001
$ws.Range("N71").Value = -19654.5
$ws.Range("H123").Value = 5009.6665
$ws.Range("I123").Value = 5009.6665
$ws.Range("K123").Value = 15028.9995
$ws.Range("M123").Value = -12578.9995

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 21904
$ws.Range("J33").Value = 21380
$ws.Range("L33").Value = 21380
$ws.Range("N33").Value = -21884
$ws.Range("H40").Value = 15999.5
$ws.Range("J40").Value = 15999.5
$ws.Range("L40").Value = 15999.5
$ws.Range("N40").Value = -16301.5
$ws.Range("H44").Value = 16666.334
$ws.Range("J44").Value = 16666.334
$ws.Range("L44").Value = 16666.334
$ws.Range("N44").Value = -17858.334
$ws.Range("H70").Value = 284106.5
$ws.Range("I70").Value = 284106.5
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 284106.5
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -283836.5
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 284106.5
$ws.Range("I73").Value = 284106.5
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 284106.5
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -283170.5
$ws.Range("N73").ClearContents()
$ws.Range("H102").Value = 4124.136
$ws.Range("I102").Value = 4101.8237
$ws.Range("K102").Value = 4101.8237
$ws.Range("M102").Value = -2479.8237
$ws.Range("H107").Value = 244.5
$ws.Range("I107").Value = 224.4
$ws.Range("J107").Value = 345
$ws.Range("K107").Value = 224.4
$ws.Range("L107").Value = 345
$ws.Range("M107").Value = 1695.6
$ws.Range("N107").Value = -4185
$ws.Range("H122").Value = 3973.7646
$ws.Range("I122").Value = 3313.8823
$ws.Range("J122").Value = 4633.647
$ws.Range("K122").Value = 9941.6469
$ws.Range("L122").Value = 13900.941
$ws.Range("M122").Value = -7491.6469
$ws.Range("N122").Value = -18800.941
$ws.Range("H132").Value = 6951.514
$ws.Range("I132").Value = 6732.5356
$ws.Range("J132").Value = 7827.4287
$ws.Range("K132").Value = 20197.6068
$ws.Range("L132").Value = 23482.2861
$ws.Range("M132").Value = -17667.6068
$ws.Range("N132").Value = -28542.2861

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 50000
$ws.Range("J36").Value = 50000
$ws.Range("L36").Value = 50000
$ws.Range("N36").Value = -51124
$ws.Range("H46").Value = 3865.9395
$ws.Range("I46").Value = 698.2
$ws.Range("K46").Value = 698.2
$ws.Range("M46").Value = -510.2
$ws.Range("H50").Value = 42495
$ws.Range("J50").Value = 42495
$ws.Range("L50").Value = 42495
$ws.Range("N50").Value = -43769
$ws.Range("H61").Value = 1001.6667
$ws.Range("J61").Value = 1005
$ws.Range("L61").Value = 1005
$ws.Range("N61").Value = -1409
$ws.Range("H68").Value = 2798.9092
$ws.Range("I68").Value = 2798.9092
$ws.Range("K68").Value = 2798.9092
$ws.Range("M68").Value = -2049.9092
$ws.Range("H71").Value = 2798.9092
$ws.Range("I71").Value = 2798.9092
$ws.Range("K71").Value = 13994.546
$ws.Range("M71").Value = -10250.546
$ws.Range("H113").Value = 1001.6667
$ws.Range("J113").Value = 1005
$ws.Range("L113").Value = 1005
$ws.Range("N113").Value = -5345
$ws.Range("H122").Value = 4585.933
$ws.Range("I122").Value = 4190.8335
$ws.Range("K122").Value = 12572.5005
$ws.Range("M122").Value = -10122.5005

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 44899.9
$ws.Range("I81").Value = 54750
$ws.Range("J81").Value = 38333.168
$ws.Range("K81").Value = 109500
$ws.Range("L81").Value = 76666.336
$ws.Range("M81").Value = -108439
$ws.Range("N81").Value = -78788.336
$ws.Range("H84").Value = 44899.9
$ws.Range("I84").Value = 54750
$ws.Range("J84").Value = 38333.168
$ws.Range("K84").Value = 547500
$ws.Range("L84").Value = 383331.68
$ws.Range("M84").Value = -542196
$ws.Range("N84").Value = -393939.68
$ws.Range("H126").Value = 4939.4546
$ws.Range("I126").Value = 4599.8887
$ws.Range("K126").Value = 13799.6661
$ws.Range("M126").Value = -11329.6661
$ws.Range("H132").Value = 4562.07
$ws.Range("I132").Value = 4042.7097
$ws.Range("K132").Value = 12128.1291
$ws.Range("M132").Value = -9598.1291

# ---- Header restyle: remove bold + thin border + center/top alignment from row 1 on every sheet ----
foreach ($ws in $wb.Worksheets) {
    $ws.Range("A1:N1").Style = "Normal"
}

Write-Host "edit complete"
